# Add a new row (61) to Sheet1 describing the "Open Access an der Leibniz
# Universität Hannover" dashboard, wire up its two hyperlinks (NAME-column
# style URL + ROR URL), and move the viewport/selection the author left
# behind after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- plain (non-hyperlink) cells -------------------------------------------------
# (values are entered in the same order the original author did, so the
# shared-string table grows with the same new-entry ordering)
$ws.Range("A61").Value = "Open Access an der Leibniz Universität Hannover"
$ws.Range("C61").Value = "2019–2024"
$ws.Range("D61").Value = "Leibniz University Hannover"

# --- hyperlinked cells (B61, E61) --------------------------------------------------
# Register the hyperlink relationships in column order (B then E) so the
# relationship ids line up with the author's file, but poke the cell text in
# the order the author actually typed it (E before B) so the shared-string
# table grows with the same new-entry ordering too.
$ws.Hyperlinks.Add($ws.Range("B61"), "https://www.uni-hannover.de/en/universitaet/profil/leitbild-und-strategien/forschung/open-science/open-access")
$ws.Hyperlinks.Add($ws.Range("E61"), "https://ror.org/0304hq317")

$ws.Range("E61").Value = "https://ror.org/0304hq317"
$ws.Range("E61").Style = "Link"
$ws.Range("E61").VerticalAlignment = -4160

$ws.Range("B61").Value = "https://www.uni-hannover.de/en/universitaet/profil/leitbild-und-strategien/forschung/open-science/open-access"
$ws.Range("B61").Style = "Link"
$ws.Range("B61").VerticalAlignment = -4160

$ws.Range("F61").Value = "RPO"
$ws.Range("G61").Value = "research institution(s)"
$ws.Range("H61").Value = "DEU"
$ws.Range("I61").Value = "publications"
$ws.Range("J61").Value = "N/A"
$ws.Range("K61").Value = "N/A"
$ws.Range("L61").Value = "open"
$ws.Range("M61").Value = "N/A"
$ws.Range("N61").Value = "N/A"

# --- remaining row-61 cells keep the sheet's plain top-aligned look --------------
$ws.Range("A61").VerticalAlignment = -4160
$ws.Range("C61").VerticalAlignment = -4160
$ws.Range("D61").VerticalAlignment = -4160
$ws.Range("F61:N61").VerticalAlignment = -4160

# --- viewport / selection, matching the author's last recorded state ------------
$excel.Goto($ws.Range("D22"), $true)
$ws.Range("O61").Select()
